$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.233789666666667
$ws.Range("H2").Value = 3.701369
$ws.Range("I2").Value = 0.0001664233864291757
$ws.Range("J2").Value = 0.0001664233864291757
$ws.Range("M2").Value = 0.140567
$ws.Range("N2").Value = 0.421701
$ws.Range("O2").Value = 0.07810038533383065
$ws.Range("P2").Value = 0.07810038533383065
$ws.Range("Q2").Value = 0.1734301120743333
$ws.Range("R2").Value = 1.560871008669
$ws.Range("S2").Value = 0.00001299773060867963
$ws.Range("T2").Value = 0.00001299773060867963

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.233789666666667
$ws.Range("H3").Value = 3.701369
$ws.Range("I3").Value = 0.0001664233864291757
$ws.Range("J3").Value = 0.0001664233864291757
$ws.Range("O3").Value = 0.02984383293631935
$ws.Range("P3").Value = 0.02984383293631935
$ws.Range("Q3").Value = 0.06627136689211112
$ws.Range("R3").Value = 0.5964423020290001
$ws.Range("S3").Value = 0.000004966711741288837
$ws.Range("T3").Value = 0.000004966711741288837

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.233789666666667
$ws.Range("H4").Value = 3.701369
$ws.Range("I4").Value = 0.0001664233864291757
$ws.Range("J4").Value = 0.0001664233864291757
$ws.Range("O4").Value = 0.8920557817298499
$ws.Range("P4").Value = 0.8920557817298499
$ws.Range("Q4").Value = 1.980903596578667
$ws.Range("R4").Value = 17.828132369208
$ws.Range("S4").Value = 0.0001484589440792072
$ws.Range("T4").Value = 0.0001484589440792072

# Row 5
$ws.Range("I5").Value = 0.9827534361704352
$ws.Range("J5").Value = 0.9827534361704352
$ws.Range("M5").Value = 0.140567
$ws.Range("N5").Value = 0.421701
$ws.Range("O5").Value = 0.07810038533383065
$ws.Range("P5").Value = 0.07810038533383065
$ws.Range("Q5").Value = 1024.129133732103
$ws.Range("R5").Value = 9217.162203588927
$ws.Range("S5").Value = 0.07675342205305713
$ws.Range("T5").Value = 0.07675342205305713

# Row 6
$ws.Range("I6").Value = 0.9827534361704352
$ws.Range("J6").Value = 0.9827534361704352
$ws.Range("O6").Value = 0.02984383293631935
$ws.Range("P6").Value = 0.02984383293631935
$ws.Range("S6").Value = 0.02932912936666425
$ws.Range("T6").Value = 0.02932912936666425

# Row 7
$ws.Range("I7").Value = 0.9827534361704352
$ws.Range("J7").Value = 0.9827534361704352
$ws.Range("O7").Value = 0.8920557817298499
$ws.Range("P7").Value = 0.8920557817298499
$ws.Range("S7").Value = 0.8766708847507138
$ws.Range("T7").Value = 0.8766708847507138

# Row 8
$ws.Range("I8").Value = 0.01708014044313564
$ws.Range("J8").Value = 0.01708014044313564
$ws.Range("M8").Value = 0.140567
$ws.Range("N8").Value = 0.421701
$ws.Range("O8").Value = 0.07810038533383065
$ws.Range("P8").Value = 0.07810038533383065
$ws.Range("Q8").Value = 17.79924525546767
$ws.Range("R8").Value = 160.193207299209
$ws.Range("S8").Value = 0.001333965550164839
$ws.Range("T8").Value = 0.001333965550164839

# Row 9
$ws.Range("I9").Value = 0.01708014044313564
$ws.Range("J9").Value = 0.01708014044313564
$ws.Range("O9").Value = 0.02984383293631935
$ws.Range("P9").Value = 0.02984383293631935
$ws.Range("S9").Value = 0.0005097368579138117
$ws.Range("T9").Value = 0.0005097368579138117

# Row 10
$ws.Range("I10").Value = 0.01708014044313564
$ws.Range("J10").Value = 0.01708014044313564
$ws.Range("O10").Value = 0.8920557817298499
$ws.Range("P10").Value = 0.8920557817298499
$ws.Range("S10").Value = 0.01523643803505699
$ws.Range("T10").Value = 0.01523643803505699
